$d = $word.ActiveDocument

$pairs = @(
    @("635÷4=", "454÷2="),
    @("437÷4=", "416÷3="),
    @("884÷5=", "560÷3="),
    @("594÷7=", "971÷3="),
    @("719÷5=", "167÷2="),
    @("562÷3=", "747÷3="),
    @("968÷8=", "371÷4="),
    @("266÷3=", "520÷7="),
    @("412÷7=", "603÷9="),
    @("268÷4=", "558÷8="),
    @("485÷5=", "152÷7="),
    @("412÷3=", "878÷6="),
    @("230÷7=", "607÷4="),
    @("334÷9=", "985÷8="),
    @("205÷3=", "286÷4="),
    @("170÷5=", "733÷3="),
    @("756÷8=", "751÷2="),
    @("975÷7=", "791÷8="),
    @("144÷3=", "702÷6="),
    @("167÷6=", "908÷9="),
    @("925÷4=", "480÷5="),
    @("107÷6=", "654÷9="),
    @("856÷2=", "309÷4="),
    @("373÷3=", "174÷2="),
    @("867÷3=", "598÷2=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
